$d = $word.ActiveDocument

# --- Change 1: ProductVariant row -----------------------------------------
# "sku, " -> "sku, default, "  (adds a new "default" attribute before product_id)
$d.Content.Find.Execute("sku, ", $true, $false, $false, $false, $false, $true, 1, $false, "sku, default, ", 2) | Out-Null

# --- Change 2: Service row --------------------------------------------------
# ", service_description" stays the same text but the two runs that produced
# it (", " and "service_description") collapse into a single run.
$d.Content.Find.Execute(", service_description", $true, $false, $false, $false, $false, $true, 1, $false, ", service_description", 2) | Out-Null
